$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 17-29 are being re-ordered ("refactored") and a couple of the D-column
# multiplicities change from 2 to 1.  Rewrite each of the 13 rows explicitly
# so that the final layout matches the target state.
# ---------------------------------------------------------------------------

# Row 17 (was row 20): network_df / 363
$ws.Range("A17").Value = 363
$ws.Range("B17").Value = "network_df"
$ws.Range("C17").Formula = "=(B11*B3*2 +B11*B6)/B2"
$ws.Range("D17").Value = 1
$ws.Range("E17").Formula = "=C17*D17"

# Row 18 (was row 21): node_1_names / 364
$ws.Range("A18").Value = 364
$ws.Range("B18").Value = "node_1_names"
$ws.Range("C18").Formula = "=B12*B3/B2"
$ws.Range("D18").Value = 1
$ws.Range("E18").Formula = "=C18*D18"

# Row 19 (was row 22): node_2_names / 364
$ws.Range("A19").Value = 364
$ws.Range("B19").Value = "node_2_names"
$ws.Range("C19").Formula = "=B12*B3/B2"
$ws.Range("D19").Value = 1
$ws.Range("E19").Formula = "=C19*D19"

# Row 20 (was row 23): unique_gene_names / 365
$ws.Range("A20").Value = 365
$ws.Range("B20").Value = "unique_gene_names"
$ws.Range("C20").Formula = "=B12*B3/B2"
$ws.Range("D20").Value = 1
$ws.Range("E20").Formula = "=C20*D20"

# Row 21 (was row 24): genes_lookup_table / 369
$ws.Range("A21").Value = 369
$ws.Range("B21").Value = "genes_lookup_table"
$ws.Range("C21").Formula = "=(B12*C12*B4)/B2"
$ws.Range("D21").Value = 1
$ws.Range("E21").Formula = "=C21*D21"

# Row 22 (was row 25): network_df / 371
$ws.Range("A22").Value = 371
$ws.Range("B22").Value = "network_df"
$ws.Range("C22").Formula = "=C17"
$ws.Range("D22").Value = 1
$ws.Range("E22").Formula = "=C22*D22"

# Row 23 (was row 26): network_df / 372
$ws.Range("A23").Value = 372
$ws.Range("B23").Value = "network_df"
$ws.Range("C23").Formula = "=C17"
$ws.Range("D23").Value = 1
$ws.Range("E23").Formula = "=C23*D23"

# Row 24 (was row 27): network_df / 374
$ws.Range("A24").Value = 374
$ws.Range("B24").Value = "network_df"
$ws.Range("C24").Formula = "=C17*2"
$ws.Range("D24").Value = 1
$ws.Range("E24").Formula = "=C24*D24"

# Row 25 (was row 28): network_mat_sparse / 374
$ws.Range("A25").Value = 374
$ws.Range("B25").Value = "network_mat_sparse"
$ws.Range("C25").Formula = "=(B12*B7+B11*B7+B11*B6)/B2"
$ws.Range("D25").Value = 1
$ws.Range("E25").Formula = "=C25*D25"

# Row 26 (was row 29): network_mat / 378
$ws.Range("A26").Value = 378
$ws.Range("B26").Value = "network_mat"
$ws.Range("C26").Formula = "=C25"
$ws.Range("D26").Value = 1
$ws.Range("E26").Formula = "=C26*D26"

# Row 27 (was row 17): drug_response_df / 359
$ws.Range("A27").Value = 359
$ws.Range("B27").Value = "drug_response_df"
$ws.Range("C27").Formula = "=(B9*B8*B6+B9*B3+B8*B3)/B2"
$ws.Range("D27").Value = 1
$ws.Range("E27").Formula = "=C27*D27"

# Row 28 (was row 18): spreadsheet_df / 360 (multiplicity now 1, was 2)
$ws.Range("A28").Value = 360
$ws.Range("B28").Value = "spreadsheet_df"
$ws.Range("C28").Formula = "=(B10*B9*B6+B10*B3+B9*B3)/B2"
$ws.Range("D28").Value = 1
$ws.Range("E28").Formula = "=C28*D28"

# Row 29 (was row 19): spreadsheet_genes_as_input / 361
$ws.Range("A29").Value = 361
$ws.Range("B29").Value = "spreadsheet_genes_as_input"
$ws.Range("C29").Formula = "=B10*B3/B2"
$ws.Range("D29").Value = 1
$ws.Range("E29").Formula = "=C29*D29"

# ---------------------------------------------------------------------------
# Downstream formulas that referenced the old row numbers must follow the
# moved cells to their new homes.
# ---------------------------------------------------------------------------
$ws.Range("C38").Formula = "=C34+C33+C32+C31+C26+C29"
$ws.Range("C40").Formula = "=C26"
$ws.Range("C71").Formula = "=C29+C20"

# ---------------------------------------------------------------------------
# Selection / view bookkeeping to match the author's saved state.
# ---------------------------------------------------------------------------
$ws.Range("B38").Select()
